# radx-rad-metadata-compiler template: make the "creator_profileid" ->
# RADx Metadata Specification Path row (row 21) return a specific type,
# i.e. point it at "/Data File Creators/Creator Identifier" instead of the
# generic "/Auxiliary Metadata/Data File Descriptive Key-Value Pairs"
# fallback, and give the row the same highlighted look as the other
# "identifier" field rows (e.g. row 6) by reusing their existing cell
# format (fill/border) rather than inventing a new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (A6/B6) already carries the target look (theme fill + thin left/
# top/bottom border, no right border) that row 21 should adopt. Copy its
# format onto A21 and B21 so the workbook's existing style table is reused
# instead of a duplicate being appended.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# A21 keeps its label ("creator_profileid"); only B21's mapped path
# changes, to the Creator Identifier path (specific type) instead of the
# generic Key-Value Pairs fallback.
$ws.Range("B21").Value = "/Data File Creators/Creator Identifier"

# Reflect the edited row as the active selection, like it was in the
# authored workbook.
$ws.Range("A21:B21").Select() | Out-Null
